$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value (45181) for every data
# row. The update bumps that value by one day (45181 -> 45182) for every
# row that currently has it, leaving everything else untouched.

$lastRow = $ws.Cells.Item(1, 1).End(4).Row
if ($lastRow -lt 2) { $lastRow = 495 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
